$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$reqLOM3234 = "LOM3234 -  Óptica Física  (Requisito)`n"
$reqLOM3259 = "LOM3259 -  Materiais e Dispositivos Eletrônicos  (Indicação de Conjunto)`n"

# Swap order: LOM3234 now comes first (row 24), LOM3259 second (row 25)
$ws.Range("B24").Value = $reqLOM3234
$ws.Range("C24").Value = $reqLOM3234

$ws.Range("B25").Value = $reqLOM3259
$ws.Range("C25").Value = $reqLOM3259
